# Update the "cryptos" worksheet with refreshed Price (column D) and
# Volume(1h) (column E) figures, as published by the scheduled GitHub
# Actions scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $cell = $ws.Range($CellRef)
    # Some of the refreshed "Price" strings (e.g. "528.43", "1.00") are
    # syntactically valid numbers/dates. Excel's COM Value setter would
    # normally auto-convert such text into a numeric/date cell, which
    # would change the cell's underlying type away from the plain text
    # that the source data represents. Forcing a text number format
    # while assigning the value keeps it as text, and re-applying the
    # "Normal" style immediately afterwards removes the temporary
    # formatting override again, leaving the cell's style untouched.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# --- Row 2: Bitcoin ---
Set-TextValue "D2" '57.630.71'
$ws.Range("E2").Value = '  -1.38%  '

# --- Row 3: Ethereum ---
Set-TextValue "D3" '3.101.99'
$ws.Range("E3").Value = '  -2.00%  '

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = '  -0.06%  '

# --- Row 5: BNB ---
Set-TextValue "D5" '528.43'
$ws.Range("E5").Value = '  -1.01%  '

# --- Row 6: Solana ---
Set-TextValue "D6" '137.60'
$ws.Range("E6").Value = '  -3.33%  '

# --- Row 7: USDC ---
Set-TextValue "D7" '1.00'
$ws.Range("E7").Value = '  -0.11%  '

# --- Row 8: LidoStakedEther ---
Set-TextValue "D8" '3.104.56'
$ws.Range("E8").Value = '  -2.05%  '

# --- Row 9: XRP ---
Set-TextValue "D9" '0.467'
$ws.Range("E9").Value = '  +4.18%  '

# --- Row 10: Toncoin ---
$ws.Range("E10").Value = '  +0.48%  '

# --- Row 11: Dogecoin ---
$ws.Range("E11").Value = '  -3.01%  '

# --- Row 12: Cardano ---
Set-TextValue "D12" '0.408'
$ws.Range("E12").Value = '  +2.10%  '

# --- Row 13: TRON ---
$ws.Range("E13").Value = '  +1.87%  '

# --- Row 14: WrappedliquidstakedEther2.0 ---
Set-TextValue "D14" '3.637.57'
$ws.Range("E14").Value = '  -2.18%  '

# --- Row 15: Avalanche ---
$ws.Range("E15").Value = '  -1.59%  '

# --- Row 16: ShibaInu ---
$ws.Range("E16").Value = '  -3.01%  '

# --- Row 17: WrappedBTC ---
Set-TextValue "D17" '57.681.13'
$ws.Range("E17").Value = '  -1.55%  '

# --- Row 18: WrappedEther ---
Set-TextValue "D18" '3.086.49'
$ws.Range("E18").Value = '  -2.73%  '

# --- Row 19: Polkadot ---
$ws.Range("E19").Value = '  -3.75%  '

# --- Row 20: Chainlink ---
$ws.Range("E20").Value = '  -2.81%  '

# --- Row 21: Uniswap ---
$ws.Range("E21").Value = '  -1.78%  '

# --- Row 22: BitcoinCash ---
Set-TextValue "D22" '350.16'
$ws.Range("E22").Value = '  -2.17%  '

# --- Row 23: Dai ---
$ws.Range("E23").Value = '  +0.09%  '

# --- Row 24: Litecoin ---
Set-TextValue "D24" '68.89'
$ws.Range("E24").Value = '  +0.82%  '

# --- Row 25: Polygon ---
$ws.Range("E25").Value = '  -2.44%  '

# --- Row 26: Kaspa ---
Set-TextValue "D26" '0.166'
$ws.Range("E26").Value = '  -2.01%  '

# --- Row 27: Binance-PegBSC-USD ---
Set-TextValue "D27" '0.999'
$ws.Range("E27").Value = '  +0.06%  '

# --- Row 28: PEPE ---
$ws.Range("E28").Value = '  -9.33%  '

# --- Row 29: USDe ---
$ws.Range("E29").Value = '  +0.01%  '

# --- Row 30: InternetComputer(DFINITY) ---
Set-TextValue "D30" '7.13'
$ws.Range("E30").Value = '  -5.09%  '

# --- Row 31: PancakeSwap ---
$ws.Range("E31").Value = '  -2.33%  '

# --- Row 32: RenderToken ---
Set-TextValue "D32" '5.99'
$ws.Range("E32").Value = '  -8.43%  '

# --- Row 33: EthereumClassic ---
Set-TextValue "D33" '21.10'
$ws.Range("E33").Value = '  -1.17%  '

# --- Row 34: NEARProtocol ---
Set-TextValue "D34" '4.88'
$ws.Range("E34").Value = '  -0.36%  '

# --- Row 35: Fetch.AI ---
$ws.Range("E35").Value = '  -6.73%  '

# --- Row 36: Monero ---
Set-TextValue "D36" '159.13'
$ws.Range("E36").Value = '  +0.98%  '

# --- Row 37: Aptos ---
Set-TextValue "D37" '6.03'
$ws.Range("E37").Value = '  -3.21%  '

# --- Row 38: EnergySwap ---
Set-TextValue "D38" '25.70'
$ws.Range("E38").Value = '  -2.61%  '

# --- Row 39: ImmutableX ---
$ws.Range("E39").Value = '  -4.53%  '

# --- Row 40: Stacks ---
Set-TextValue "D40" '1.63'
$ws.Range("E40").Value = '  -1.98%  '

# --- Row 41: Hedera ---
$ws.Range("E41").Value = '  -1.51%  '

# --- Row 42: Filecoin ---
$ws.Range("E42").Value = '  -0.08%  '

# --- Row 43: Mantle ---
Set-TextValue "D43" '0.694'
$ws.Range("E43").Value = '  -2.00%  '

# --- Row 44: Maker ---
Set-TextValue "D44" '2.398.64'
$ws.Range("E44").Value = '  +2.41%  '

# --- Row 45: RenzoRestakedETH ---
$ws.Range("E45").Value = '  +0.23%  '

# --- Row 46: VeChain ---
$ws.Range("E46").Value = '  +0.05%  '

# --- Row 47: Cosmos ---
Set-TextValue "D47" '3.144.59'
$ws.Range("E47").Value = '  -2.10%  '

# --- Row 48: ONDO ---
Set-TextValue "D48" '0.0262'
$ws.Range("E48").Value = '  -4.32%  '

# --- Row 49: InjectiveProtocol ---
$ws.Range("E49").Value = '  -0.83%  '

# --- Row 50: (unnamed) ---
$ws.Range("E50").Value = '  -6.12%  '

# --- Row 51: (unnamed) ---
Set-TextValue "D51" '19.56'
$ws.Range("E51").Value = '  -5.30%  '
